$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "47÷8=" "15÷2="
Replace-Text "43÷2=" "91÷9="
Replace-Text "84÷5=" "31÷4="
Replace-Text "18÷8=" "33÷5="
Replace-Text "71÷3=" "59÷9="
Replace-Text "65÷2=" "54÷9="
Replace-Text "41÷4=" "72÷9="
Replace-Text "64÷4=" "73÷4="
Replace-Text "54÷8=" "99÷7="
Replace-Text "32÷6=" "24÷9="
Replace-Text "70÷4=" "63÷8="
Replace-Text "96÷9=" "56÷9="
Replace-Text "73÷7=" "47÷9="
Replace-Text "61÷2=" "76÷8="
Replace-Text "10÷2=" "31÷8="
Replace-Text "42÷5=" "15÷5="
Replace-Text "34÷3=" "52÷6="
Replace-Text "21÷9=" "99÷8="
Replace-Text "37÷7=" "84÷8="
Replace-Text "90÷5=" "43÷3="
Replace-Text "60÷3=" "12÷9="
Replace-Text "50÷2=" "17÷2="
Replace-Text "62÷6=" "54÷9="
Replace-Text "60÷7=" "60÷6="
Replace-Text "52÷2=" "19÷9="
